$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update while forcing text storage (matches the
# source workbook, where these cells are plain strings, not numbers),
# then restore the original cell style so no formatting is disturbed.

$cell = $ws.Cells.Item(2, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '43.040.37'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(2, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.99%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.548.15'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.03%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(4, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(4, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.15%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '299.61'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.86%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '92.72'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(6, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -6.04%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(7, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.25%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(8, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(9, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.62%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '35.80'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(10, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -7.13%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0806'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(11, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.78%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.68'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.41%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(13, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.23%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.932.83'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.524.44'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.62%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(16, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.23%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '14.11'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.33%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(18, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '43.023.52'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(18, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.90%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.12'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.65%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0984'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(20, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.65%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.61'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(21, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.91%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '71.81'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(22, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.32%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '257.27'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -9.64%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(24, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.91%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(25, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.99%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '29.15'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(27, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.04'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(28, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.94%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '37.48'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.53%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(30, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -5.76%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.96'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(31, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.28%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '153.95'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.02%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 2)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 3)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.17'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.81%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 2)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'WEMIXToken'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 3)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.75'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.31%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.38'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -6.25%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0801'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(36, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.44%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(37, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.95%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(38, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.17%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.03'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(39, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +7.66%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '23.48'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(40, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +9.62%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.45'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(41, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.11%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.91'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.25%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0311'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(43, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -4.19%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.081.78'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.78%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '84.68'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(46, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -9.95%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.90'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.08%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(48, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.81%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(49, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.791.13'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(49, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -1.89%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '104.85'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.34%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -3.62%  '
$cell.Style = $origStyle
